# Append a new "1.1" version row to the Versions table on Sheet1.
#
# Target state (per the commit "general: - add the icon - report on the
# failed contacts in the summary - trim the template column while
# generating - reset the lookups on Go - tell the user that the csv file
# is open"):
#   - Table1 / the sheet's used range grows from A1:C2 to A1:C3
#   - A3 = 1.1 (stored as a number, formatted as Text like A2)
#   - B3 = the multi-line release-notes string (wrapped, row tall enough
#     for all 6 lines)
#   - C3 stays empty (not yet released)
#   - selection ends up on B3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the table by one row so Table1's ref / autoFilter / dimension all
# update together, instead of poking sheetData directly.
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

# Version number, formatted as text (same numFmt "@" style as A2's 1.0).
$ws.Cells.Item(3, 1).Value = 1.1
$ws.Cells.Item(3, 1).NumberFormat = "@"

# Release notes, wrapped across 6 lines.
$notes = "general:`n- add the icon`n- report on the failed contacts in the summary`n- trim the template column while generating`n- reset the lookups on Go`n- tell the user that the csv file is open"
$ws.Cells.Item(3, 2).Value = $notes
$ws.Cells.Item(3, 2).WrapText = $true
$ws.Cells.Item(3, 2).EntireRow.RowHeight = 90

# Match the saved selection.
$ws.Range("B3").Select() | Out-Null
